# Commit: "renaming all caps variables to lowercase to match all other variables"
#
# The ENVELOPE.xlsx workbook has several sheets (WINDOW, ROOF, WALL, FLOOR)
# that each carry a pair of header cells named GHG_<TYPE>_kgCO2m2 /
# Service_Life_<TYPE> where <TYPE> was written in ALL CAPS
# (WIN/ROOF/WALL/FLOOR). This rewrites those headers to lowercase
# (win/roof/wall/floor) to match the naming convention used everywhere
# else in the workbook, and leaves the active sheet/selection on SHADING
# (matching the editor's on-save cursor state).

$wb = $excel.ActiveWorkbook

$wsWindow = $wb.Worksheets.Item("WINDOW")
$wsRoof   = $wb.Worksheets.Item("ROOF")
$wsWall   = $wb.Worksheets.Item("WALL")
$wsFloor  = $wb.Worksheets.Item("FLOOR")
$wsShading = $wb.Worksheets.Item("SHADING")

# WALL: G1/H1 headers
$wsWall.Range("G1").Value = "GHG_wall_kgCO2m2"
$wsWall.Range("H1").Value = "Service_Life_wall"

# FLOOR: D1/E1 headers
$wsFloor.Range("D1").Value = "GHG_floor_kgCO2m2"
$wsFloor.Range("E1").Value = "Service_Life_floor"

# WINDOW: G1/H1 headers
$wsWindow.Range("G1").Value = "GHG_win_kgCO2m2"
$wsWindow.Range("H1").Value = "Service_Life_win"

# ROOF: G1/H1 headers
$wsRoof.Range("G1").Value = "GHG_roof_kgCO2m2"
$wsRoof.Range("H1").Value = "Service_Life_roof"

# Restore the selections left on each sheet after the edits
$null = $wsWindow.Range("G2").Select()
$null = $wsRoof.Range("G2").Select()
$null = $wsWall.Range("B2").Select()
$null = $wsFloor.Range("D2").Select()

# SHADING ends up the active tab with D21 selected
$null = $wsShading.Activate()
$null = $wsShading.Range("D21").Select()
